# Commit: Added Extent report with thread safe and error screenshot
#
# Updates the TestCases workbook:
#  - TestSuite: stop re-running TC02/TC03/TC04 (Run column -> "No"),
#    leave TC01 as "Yes"
#  - TC01: "Verify home Page Header" step value "First Name" -> "FirstName"
#  - TC02: locator changed to use single quotes in the xpath
#  - Sheet selections / active tab updated to reflect where the author
#    left off (TC02 active)

$wb = $excel.ActiveWorkbook

# --- TestSuite sheet: flip Run flag to "No" for TC02, TC03, TC04 ---
$testSuite = $wb.Worksheets.Item("TestSuite")
$testSuite.Range("C3").Value = "No"
$testSuite.Range("C4").Value = "No"
$testSuite.Range("C5").Value = "No"
$testSuite.Range("C6").Select()

# --- TC01 sheet: Value column for the header-verification step ---
$tc01 = $wb.Worksheets.Item("TC01")
$tc01.Activate()
$tc01.Range("E14").Value = "FirstName"
$tc01.Range("E16").Select()

# --- TC02 sheet: LocatorValue now uses single quotes ---
$tc02 = $wb.Worksheets.Item("TC02")
$tc02.Activate()
$tc02.Range("C14").Value = "//label[text()='First Name']"
$tc02.Range("E10").Select()

# --- TC03 sheet: selection moved ---
$tc03 = $wb.Worksheets.Item("TC03")
$tc03.Activate()
$tc03.Range("E16").Select()

# Leave TC02 as the active/selected tab, matching the saved workbook view
$tc02.Activate()
